$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns permuted across rows 2-36: D(4), M(13), N(14), O(15), P(16), R(18), S(19)
$cols = @(4, 13, 14, 15, 16, 18, 19)

# Snapshot current values for rows 2-36 before writing, since this is a full-row permutation
$snapshot = @{}
for ($r = 2; $r -le 36; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# Destination row -> source row (values now located at dest row came from source row)
$mapping = @{
    2 = 34
    3 = 14
    4 = 5
    5 = 25
    6 = 31
    7 = 4
    8 = 12
    9 = 18
    10 = 30
    11 = 28
    12 = 35
    13 = 8
    14 = 19
    15 = 9
    16 = 20
    17 = 15
    18 = 3
    19 = 32
    20 = 33
    21 = 7
    22 = 29
    23 = 26
    24 = 17
    25 = 24
    26 = 36
    27 = 16
    28 = 11
    29 = 2
    30 = 27
    31 = 13
    32 = 23
    33 = 21
    34 = 22
    35 = 6
    36 = 10
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value2 = $srcVals[$c]
    }
}
